# Update column G ("K") values for rows 2-15 on Sheet1
# per the commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 3
    6  = 1
    7  = 0
    8  = 1
    9  = 1
    10 = 3
    11 = 2
    12 = 2
    13 = 1
    14 = 0
    15 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
